# Update cryptocurrency price/volume figures per the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "317.73"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-3.27%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "41.96"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-4.77%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.206"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-3.67%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08134"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-2.96%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.369"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-1.43%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.753"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-10.16%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9318"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-4.50%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1122"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-1.02%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1856"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-2.32%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09318"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-4.80%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.04580"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-1.68%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.410"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-19.17%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.1055"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.78%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001314"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.71%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005920"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-3.10%"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-1.17%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3375"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "1.33%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.95%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "2.02%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.04169"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "0.13%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.001243"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-3.98%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.004272"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-3.05%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0001223"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-5.94%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0002985"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.05%"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02588"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "-2.74%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05492"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-2.78%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.008043"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "2.66%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1395"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-1.56%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.006512"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-11.49%"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-0.87%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008233"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "4.26%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3476"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-1.15%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006737"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-1.24%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000752"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.26%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003394"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-3.29%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.004109"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "16.45%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002105"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.26%"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.26%"
